# bench on ss512 & ss1024
# Adds two new benchmark blocks (curve SS512, curve SS1024) below the
# existing secp521r1 block, and normalises a couple of leftover cell
# styles that are no longer distinct from the "red header" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Strip the "applyNumberFormat"-only styling left on B3:D3 / D21 ---
# B3:D3 used a style that was font(red)+applyNumberFormat; the number
# format was a no-op, so it collapses onto the plain red-font style (the
# same one already used by A13, D1, D15, D19, B17, C17, D17 ...).
$ws.Range("B3:D3").Font.Color = 255

# D21 used a style that was default-font+applyNumberFormat; the number
# format was a no-op there too, so it collapses onto the default style
# (no explicit style at all).
$d21 = $ws.Range("D21").Value2
$ws.Range("D21").ClearFormats()
$ws.Range("D21").Value2 = $d21

# --- 2. New benchmark block: curve SS512 (rows 24-26) ---
$ws.Range("A24").Value2 = "Based on CDH assumption"
$ws.Range("B24").Value2 = "curve SS512"
$ws.Range("C24").Value2 = "keygen required"
$ws.Range("D24").Value2 = 104.23850059509201
$ws.Range("D24").Font.Color = 255

$ws.Range("A25").Value2 = "Generazione 1000 messaggi (da 0 a 1000000)"
$ws.Range("B25").Value2 = "Commitment time"
$ws.Range("C25").Value2 = "Proof time"
$ws.Range("D25").Value2 = "Verifiy time"

$ws.Range("B26").Value2 = 0.108036279678344
$ws.Range("C26").Value2 = 0.106045007705688
$ws.Range("D26").Value2 = 0.0015685558319091699
$ws.Range("B26:D26").Font.Color = 255

# --- 3. New benchmark block: curve SS1024 (rows 28-30) ---
$ws.Range("A28").Value2 = "Based on CDH assumption"
$ws.Range("B28").Value2 = "curve SS1024"
$ws.Range("C28").Value2 = "keygen required"
$ws.Range("D28").Value2 = 201.82421469688401
$ws.Range("D28").Font.Color = 255

$ws.Range("A29").Value2 = "Generazione 1000 messaggi (da 0 a 1000000)"
$ws.Range("B29").Value2 = "Commitment time"
$ws.Range("C29").Value2 = "Proof time"
$ws.Range("D29").Value2 = "Verifiy time"

$ws.Range("B30").Value2 = 0.23519372940063399
$ws.Range("C30").Value2 = 0.23816847801208399
$ws.Range("D30").Value2 = 0.033162832260131801
$ws.Range("B30:D30").Font.Color = 255

# --- 4. View state: zoom + selection on the new last cell ---
$excel.ActiveWindow.Zoom = 102
$null = $ws.Range("D31").Select()
